$wb = $excel.ActiveWorkbook

# Remove the extra (empty) worksheets, keeping only Sheet1
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null
$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("Sheet1")

# Header / train info row
$ws.Range("A1").Value = "Train Information"
$ws.Range("B1").Value = "From:"
$ws.Range("C1").Value = "Bucuresti Nord"
$ws.Range("D1").Value = "To:"
$ws.Range("E1").Value = "Constanta"

# Column headers
$ws.Range("A3").Value = "Departure Time"
$ws.Range("B3").Value = "Delay"
$ws.Range("C3").Value = "Arrival Time"
$ws.Range("D3").Value = "Length"

# Schedule rows: Departure time, Delay label, Arrival time
$schedule = @(
    @{ Row = 4;  Dep = 0.26666666666666666; Delay = "+19 min întârziere";                          Arr = 0.47291666666666665 },
    @{ Row = 5;  Dep = 0.29166666666666669; Delay = "la timp";                                       Arr = 0.39513888888888887 },
    @{ Row = 6;  Dep = 0.34722222222222227; Delay = "la timp";                                       Arr = 0.44930555555555557 },
    @{ Row = 7;  Dep = 0.39583333333333331; Delay = "-1 min mai devreme";                             Arr = 0.48402777777777778 },
    @{ Row = 8;  Dep = 0.4375;              Delay = "+4 min întârziere";                              Arr = 0.5395833333333333 },
    @{ Row = 9;  Dep = 0.47916666666666669; Delay = "+2 min întârziere";                              Arr = 0.58750000000000002 },
    @{ Row = 10; Dep = 0.52083333333333337; Delay = "la timp";                                       Arr = 0.62916666666666665 },
    @{ Row = 11; Dep = 0.59722222222222221; Delay = "+20 min întârziere";                             Arr = 0.70138888888888884 },
    @{ Row = 12; Dep = 0.64027777777777783; Delay = "la timp";                                       Arr = 0.74861111111111101 },
    @{ Row = 13; Dep = 0.68958333333333333; Delay = "+4 min întârziere";                              Arr = 0.79236111111111107 },
    @{ Row = 14; Dep = 0.72916666666666663; Delay = "la timp";                                       Arr = 0.81736111111111109 },
    @{ Row = 15; Dep = 0.77083333333333337; Delay = "sosește cu 102 min întârziere la Constanța*";    Arr = 0.87916666666666676 },
    @{ Row = 16; Dep = 0.8520833333333333;  Delay = "sosește la timp la Constanța*";                  Arr = 0.95833333333333337 }
)

foreach ($entry in $schedule) {
    $r = $entry.Row

    $ws.Range("A$r").Value = $entry.Dep
    $ws.Range("A$r").NumberFormat = "h:mm"

    $ws.Range("B$r").Value = $entry.Delay

    $ws.Range("C$r").Value = $entry.Arr
    $ws.Range("C$r").NumberFormat = "h:mm"
}
